# Journal de bord - add the 25/05/2016 entry as a new bulleted paragraph,
# and move the "_GoBack" bookmark from its old spot (mid-sentence in the
# 24/05/2016 entry) to the end of the new 25/05/2016 entry.

$d = $word.ActiveDocument

# --- 1. Create the new paragraph right after the existing (only) bullet ---
$p2 = $d.Paragraphs(2)
$p2.Range.InsertParagraphAfter()

# --- 2. Fill it with the first run of text ---
$p3 = $d.Paragraphs(3)
$r3 = $p3.Range
$run1Text = "25/05/2016 : Tentative de maj de moverio bt200 vers le build developpeur (besoin acces wifi), lecture de docs"
$r3.InsertAfter($run1Text)

# Remember where run 1 ends (just before the paragraph mark) - this is
# where the second run needs to be split off into its own <w:r>.
$p3 = $d.Paragraphs(3)
$r3 = $p3.Range
$splitPos = $r3.End - 1

# --- 3. Append the second run of text right after the first ---
$endOfRun1 = $d.Range($r3.End, $r3.End)
$run2Text = ". Decouverte Wikitude et Vuforia"
$endOfRun1.InsertAfter($run2Text)

# --- 4. Force runs 1 and 2 to stay as separate <w:r> elements (matching
#        the target) by briefly bookmarking the boundary between them and
#        removing the bookmark again - the bookmark insertion/removal
#        splits the run without altering any text. ---
$d.Bookmarks.Add("TempSplit", $d.Range($splitPos, $splitPos))
$d.Bookmarks("TempSplit").Delete()

# --- 5. Place the "_GoBack" bookmark at the very end of the new paragraph.
#        Adding a zero-length bookmark exactly at the end of the document
#        content is unreliable, so pad with temporary sentinel text first,
#        add the bookmark right before the padding, then delete the
#        padding - the bookmark stays anchored at the real end point. ---
$p3 = $d.Paragraphs(3)
$r3 = $p3.Range
$sentinel = "ZZZZZZZZZZ"
$padRange = $d.Range($r3.End, $r3.End)
$padRange.InsertAfter($sentinel)

$p3 = $d.Paragraphs(3)
$r3 = $p3.Range
$targetPos = $r3.End - 1 - $sentinel.Length

# Re-adding a bookmark named "_GoBack" relocates the existing one (from
# the 24/05/2016 paragraph) to this new position, so the old occurrence
# disappears automatically.
$d.Bookmarks.Add("_GoBack", $d.Range($targetPos, $targetPos))

# Remove the sentinel padding; the bookmark stays put at $targetPos.
$p3 = $d.Paragraphs(3)
$r3 = $p3.Range
$padToDelete = $d.Range($targetPos, $r3.End - 1)
$padToDelete.Delete()
